$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and volume-change (E) columns for existing rows ---

$ws.Range("D2").Value = "69.429.86"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").Value = "3.691.14"
$ws.Range("E3").Value = "  -3.04%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "686.55"
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.20"
$ws.Range("E6").Value = "  -5.91%  "

$ws.Range("D7").Value = "3.689.88"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -5.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  -8.84%  "

$ws.Range("E11").Value = "  -3.90%  "

$ws.Range("E12").Value = "  -9.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  -7.12%  "

$ws.Range("D14").Value = "4.314.36"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.52"
$ws.Range("E15").Value = "  -10.41%  "

$ws.Range("D16").Value = "3.690.37"
$ws.Range("E16").Value = "  -3.42%  "

$ws.Range("D17").Value = "69.425.61"
$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.96"
$ws.Range("E19").Value = "  -9.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  -10.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.50"
$ws.Range("E21").Value = "  -8.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("E22").Value = "  -5.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.647"
$ws.Range("E23").Value = "  -9.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.71"
$ws.Range("E24").Value = "  -4.60%  "

$ws.Range("D25").Value = "3.839.92"
$ws.Range("E25").Value = "  -3.02%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -11.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.02"
$ws.Range("E28").Value = "  -13.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  -10.42%  "

$ws.Range("E30").Value = "  -9.81%  "

$ws.Range("E31").Value = "  -12.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.66"
$ws.Range("E32").Value = "  -8.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.02"
$ws.Range("E33").Value = "  -11.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.75"
$ws.Range("E35").Value = "  -8.27%  "

$ws.Range("E36").Value = "  -6.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.20"
$ws.Range("E37").Value = "  -12.08%  "

$ws.Range("E38").Value = "  -7.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.28"
$ws.Range("E39").Value = "  -3.26%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0905"
$ws.Range("E41").Value = "  -10.31%  "

$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "373.70"

# --- Row 43/44: swap Monero and Mantle (with updated price/volume) ---
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.943"
$ws.Range("E43").Value = "  -6.75%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "166.63"
$ws.Range("E44").Value = "  +0.91%  "

# --- Row 49/50: swap FLOKI and InjectiveProtocol (with updated price/volume) ---
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "28.43"
$ws.Range("E49").Value = "  -5.51%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000275"
$ws.Range("E50").Value = "  -9.28%  "
